$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.393.11"
$ws.Range("E2").Value = "  +6.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.725.58"
$ws.Range("E3").Value = "  +4.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.84%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.39"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3738"
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.29"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3353"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.173"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07355"
$ws.Range("E11").Value = "  +3.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.009"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.366"
$ws.Range("E13").Value = "  +5.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.12"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("E15").Value = "  +7.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.733.96"
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001071"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06627"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.63"
$ws.Range("E19").Value = "  +4.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.52"
$ws.Range("E21").Value = "  +4.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.117"
$ws.Range("E22").Value = "  +3.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.84"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.478.58"
$ws.Range("E24").Value = "  +6.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.453"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.68"
$ws.Range("E26").Value = "  +4.09%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.384"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.383"
$ws.Range("E28").Value = "  +15.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.41"
$ws.Range("E29").Value = "  +4.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.927.71"
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.06"
$ws.Range("E31").Value = "  +4.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.141"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.971"
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08573"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.693"
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.71"
$ws.Range("E36").Value = "  +3.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.373"
$ws.Range("E37").Value = "  +4.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02327"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2161"
$ws.Range("E39").Value = "  +4.42%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06212"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.486"
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.223"
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6162"
$ws.Range("E43").Value = "  +3.77%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.95"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.906"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5961"
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.40"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.037"
$ws.Range("E49").Value = "  +4.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07204"
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.85"
$ws.Range("E51").Value = "  +2.64%  "
